# testplanner: Change the default template
# Insert two new header columns into the testplan table:
#   - "Testbench"   right after "Metric"              (new column E)
#   - "Coverpoints" right after "Checking Mechanism"   (new column I, after the first insert)
# Inserting whole columns shifts everything to their right (cells, styles,
# the merged title range and the custom column-width definition) just like
# using Excel's Insert > Sheet Columns from the ribbon.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Metric" is in column D, so inserting at E pushes "Intent" (and everything
# after it) one column to the right and creates a blank E2 for the new header.
$ws.Range("E:E").EntireColumn.Insert()
$ws.Range("E2").Value = "Testbench"

# After the first insert, "Checking Mechanism" sits in column H, so the second
# new column goes in at I, right before "Assignee".
$ws.Range("I:I").EntireColumn.Insert()
$ws.Range("I2").Value = "Coverpoints"
